# ---------------------------------------------------------------------------
# master_log.xlsx update — append "today's" row to the logging sheets.
#   - Cycling:   full new workout row (row 2)
#   - Running:   new HR / zone-time header columns (K1:Q1)
#   - Nutrition: new row (row 2) with just the date filled in
#   - Checkin:   new row (row 2) with just the date filled in
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newDate = 45979

# --- Cycling -----------------------------------------------------------
$cycling = $wb.Worksheets.Item("Cycling")

$cycling.Range("A2").Value = $newDate
# Setting the format twice (lower-case, then upper-case) registers both
# formatCodes in styles.xml<numFmts>, and only the final one (upper-case)
# ends up referenced by the new cell's style.
$cycling.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$cycling.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$cycling.Range("B2").Value = "Build"
$cycling.Range("C2").Value = "Cycling"
$cycling.Range("D2").Value = "Outdoors"
$cycling.Range("E2").Value = 2
$cycling.Range("F2").Value = 56
$cycling.Range("G2").Value = 22
$cycling.Range("H2").Value = 54
$cycling.Range("I2").Value = 234
$cycling.Range("J2").Value = "Zone 2"
$cycling.Range("K2").Value = "Mix"
$cycling.Range("L2").Value = 9
$cycling.Range("M2").Value = 87
$cycling.Range("N2").Value = 89
$cycling.Range("O2").Value = 370
$cycling.Range("P2").Value = 90
$cycling.Range("Q2").Value = 3
$cycling.Range("R2").Value = 1.5
$cycling.Range("S2").Value = 156
$cycling.Range("T2").Value = 124
$cycling.Range("U2").Value = 10
$cycling.Range("V2").Value = 50
$cycling.Range("W2").Value = 60
$cycling.Range("X2").Value = 0
$cycling.Range("Y2").Value = 0

# --- Running -------------------------------------------------------------
# New HR / heart-rate-zone columns appended after the existing headers.
$running = $wb.Worksheets.Item("Running")

$running.Range("K1").Value = "Max HR"
$running.Range("L1").Value = "Avg HR"
$running.Range("M1").Value = "Z1 Time (min)"
$running.Range("N1").Value = "Z2 Time (min)"
$running.Range("O1").Value = "Z3 Time (min)"
$running.Range("P1").Value = "Z4 Time (min)"
$running.Range("Q1").Value = "Z5 Time (min)"

# Match the bold/centered/bordered header style used by the rest of row 1.
$running.Range("A1").Copy()
$running.Range("K1:Q1").PasteSpecial(-4122)

# --- Nutrition -------------------------------------------------------------
$nutrition = $wb.Worksheets.Item("Nutrition")
$nutrition.Range("A2").Value = $newDate

# Reuse the date style created above instead of re-deriving it, so no
# duplicate numFmt/cellXf entries get created.
$cycling.Range("A2").Copy()
$nutrition.Range("A2").PasteSpecial(-4122)
$nutrition.Range("A2").Value = $newDate

$nutrition.Range("B2").Value = ""
$nutrition.Range("C2").Value = ""
$nutrition.Range("D2").Value = ""
$nutrition.Range("E2").Value = ""
$nutrition.Range("F2").Value = ""
$nutrition.Range("G2").Value = ""
$nutrition.Range("H2").Value = ""
$nutrition.Range("I2").Value = ""

# --- Checkin -----------------------------------------------------------
$checkin = $wb.Worksheets.Item("Checkin")
$checkin.Range("A2").Value = $newDate

$cycling.Range("A2").Copy()
$checkin.Range("A2").PasteSpecial(-4122)
$checkin.Range("A2").Value = $newDate

$checkin.Range("B2").Value = ""
$checkin.Range("C2").Value = ""
$checkin.Range("D2").Value = ""
$checkin.Range("E2").Value = ""
$checkin.Range("F2").Value = ""
$checkin.Range("G2").Value = ""
$checkin.Range("H2").Value = ""
$checkin.Range("I2").Value = ""
$checkin.Range("J2").Value = ""

Write-Host "master_log.xlsx updated"
